# "Go to bar" -> "Pay to bar" and "Go to restaurant" -> "Pay to restaurant".
# The "_GoBack" bookmark (Word's "last edit location" marker) moves from the
# top of the document to sit right after the "Pay" that replaces "Go" in
# "Go to restaurant" -- that is where the author's cursor was when they
# finished editing.

$d = $word.ActiveDocument

function Replace-GoWithPay($paraIndex) {
    $p = $d.Paragraphs.Item($paraIndex)
    $pStart = $p.Range.Start
    $text = $p.Range.Text
    $idx = $text.IndexOf("Go")
    $goStart = $pStart + $idx
    $goEnd = $goStart + 2

    # Bookmark right before "Go" keeps it from merging back into the
    # preceding run (the single space run) once it is retyped.
    $bmLeft = $d.Range($goStart, $goStart)
    $d.Bookmarks.Add("TempLeft", $bmLeft)

    $r = $d.Range($goStart, $goEnd)
    $r.Text = "Pay"

    # Return the position right after the freshly typed "Pay".
    return $goStart + 3
}

# Drop the old "_GoBack" bookmark near the top of the document; it will be
# re-added at the new edit location below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# "a. Go to bar" -> "a. Pay to bar"
$midBar = Replace-GoWithPay 73
$bmMidBar = $d.Range($midBar, $midBar)
$d.Bookmarks.Add("TempMidBar", $bmMidBar)
$d.Bookmarks.Item("TempLeft").Delete()
$d.Bookmarks.Item("TempMidBar").Delete()

# "b. Go to restaurant" -> "b. Pay to restaurant"
$midRestaurant = Replace-GoWithPay 74
$bmGoBack = $d.Range($midRestaurant, $midRestaurant)
$d.Bookmarks.Add("_GoBack", $bmGoBack)
$d.Bookmarks.Item("TempLeft").Delete()
